$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 2.93691712622856
$ws.Range("C2").Value = 2.81607473262677
$ws.Range("D2").Value = 2.76606865407216
$ws.Range("E2").Value = 2.74049145876704
$ws.Range("F2").Value = 2.71530968165729
$ws.Range("G2").Value = 2.67398272440304
$ws.Range("H2").Value = 2.6287454406253
$ws.Range("I2").Value = 2.60614002381782
$ws.Range("J2").Value = 2.59241095800179
$ws.Range("K2").Value = 2.57872082501781
$ws.Range("L2").Value = 2.56623369362173
$ws.Range("M2").Value = 2.552748715873
$ws.Range("N2").Value = 2.51556802740449
$ws.Range("O2").Value = 2.46569986202377
$ws.Range("P2").Value = 2.41756635222266
$ws.Range("Q2").Value = 2.37045340047381
$ws.Range("R2").Value = 2.32492813941118
$ws.Range("S2").Value = 2.28129139551715
$ws.Range("T2").Value = 2.23972161754873
$ws.Range("U2").Value = 2.19953712429269
$ws.Range("V2").Value = 2.16045354180412
$ws.Range("W2").Value = 2.12396132274222
$ws.Range("X2").Value = 2.08925708731948
$ws.Range("Y2").Value = 2.05696705664333
$ws.Range("Z2").Value = 2.02590362135109
$ws.Range("AA2").Value = 1.99613930838994
$ws.Range("AB2").Value = 1.9684419597018
$ws.Range("AC2").Value = 1.94184783233608
$ws.Range("AD2").Value = 0.784848519635498

# Row 3
$ws.Range("B3").Value = 0.0174107203755628
$ws.Range("C3").Value = 0.0201016509404726
$ws.Range("D3").Value = 0.0213655229741569
$ws.Range("E3").Value = 0.0208658050099303
$ws.Range("F3").Value = 0.0200756496155962
$ws.Range("G3").Value = 0.0202289682200405
$ws.Range("H3").Value = 0.020397454389265
$ws.Range("I3").Value = 0.020484119816178
$ws.Range("J3").Value = 0.0201681362919968
$ws.Range("K3").Value = 0.0199115018599077
$ws.Range("L3").Value = 0.0196727624834704
$ws.Range("M3").Value = 0.0195410359643192
$ws.Range("N3").Value = 0.0196994858524693
$ws.Range("O3").Value = 0.0200074680077551
$ws.Range("P3").Value = 0.0203782660894666
$ws.Range("Q3").Value = 0.0207864896572341
$ws.Range("R3").Value = 0.0212229678802252
$ws.Range("S3").Value = 0.0216988469269334
$ws.Range("T3").Value = 0.0221638331903511
$ws.Range("U3").Value = 0.0226570814249345
$ws.Range("V3").Value = 0.0231802589834647
$ws.Range("W3").Value = 0.0236799296660227
$ws.Range("X3").Value = 0.0241397313035659
$ws.Range("Y3").Value = 0.0245253857014175
$ws.Range("Z3").Value = 0.0248649711219498
$ws.Range("AA3").Value = 0.0251494353072223
$ws.Range("AB3").Value = 0.0253710258179548
$ws.Range("AC3").Value = 0.0254967619380203
$ws.Range("AD3").Value = 0.0623010248366342

# Row 4
$ws.Range("B4").Value = 0.0233616278523883
$ws.Range("C4").Value = 0.0261237734386893
$ws.Range("D4").Value = 0.0269161303856147
$ws.Range("E4").Value = 0.0268887089148494
$ws.Range("F4").Value = 0.0267956742431609
$ws.Range("G4").Value = 0.0271834153146973
$ws.Range("H4").Value = 0.0276364569773731
$ws.Range("I4").Value = 0.0278933249452431
$ws.Range("J4").Value = 0.0279480014514518
$ws.Range("K4").Value = 0.0280174457206488
$ws.Range("L4").Value = 0.0280665391358355
$ws.Range("M4").Value = 0.0281578173439868
$ws.Range("N4").Value = 0.0285328142310098
$ws.Range("O4").Value = 0.0290675895565852
$ws.Range("P4").Value = 0.029594828575606
$ws.Range("Q4").Value = 0.0301375709453198
$ws.Range("R4").Value = 0.0306805847083339
$ws.Range("S4").Value = 0.0312090934417946
$ws.Range("T4").Value = 0.0317099688729898
$ws.Range("U4").Value = 0.0322029322458677
$ws.Range("V4").Value = 0.0326938681732859
$ws.Range("W4").Value = 0.0331269946672791
$ws.Range("X4").Value = 0.0335301057372997
$ws.Range("Y4").Value = 0.0338708294950747
$ws.Range("Z4").Value = 0.0341889682622026
$ws.Range("AA4").Value = 0.0344760703558556
$ws.Range("AB4").Value = 0.0347004574684324
$ws.Range("AC4").Value = 0.0349007067973631

# Row 5
$ws.Range("B5").Value = 0.855977809072751
$ws.Range("C5").Value = 0.851263019296689
$ws.Range("D5").Value = 0.851184142158449
$ws.Range("E5").Value = 0.854051126013831
$ws.Range("F5").Value = 0.857170357562284
$ws.Range("G5").Value = 0.749712917055826
$ws.Range("H5").Value = 0.000000131007853471637
$ws.Range("I5").Value = 0.000000131685171253556
$ws.Range("J5").Value = 0.000000131899699335151
$ws.Range("K5").Value = 0.000000132147186848787
$ws.Range("L5").Value = 0.000000132374138481189
$ws.Range("M5").Value = 0.00000013266346502261
$ws.Range("N5").Value = 0.000000134243055491764
$ws.Range("O5").Value = 0.00000013658846803124
$ws.Range("P5").Value = 0.000000138979629946828
$ws.Range("Q5").Value = 0.000000141446122536859
$ws.Range("R5").Value = 0.00000014399058147728
$ws.Range("S5").Value = 0.0000001465980647734
$ws.Range("T5").Value = 0.0000001492514543921
$ws.Range("U5").Value = 0.000000151966372705659
$ws.Range("V5").Value = 0.000000154742024277244
$ws.Range("W5").Value = 0.000000157549009820717
$ws.Range("X5").Value = 0.000000160418343296072
$ws.Range("Y5").Value = 0.000000163318343711911
$ws.Range("Z5").Value = 0.000000166267987623986
$ws.Range("AA5").Value = 0.000000169252101982063
$ws.Range("AB5").Value = 0.00000017226088419158
$ws.Range("AC5").Value = 0.000000175325825842697
$ws.Range("AD5").Value = 0.474506366225413

# Row 6
$ws.Range("B6").Value = 0.944750157300702
$ws.Range("C6").Value = 0.94548844367585
$ws.Range("D6").Value = 0.947465795518221
$ws.Range("E6").Value = 0.949805639938611
$ws.Range("F6").Value = 0.952041681421042
$ws.Range("G6").Value = 0.845125300590563
$ws.Range("H6").Value = 0.0960340423744916
$ws.Range("I6").Value = 0.0963775764465923
$ws.Range("J6").Value = 0.0961162696431479
$ws.Range("K6").Value = 0.0959290797277433
$ws.Range("L6").Value = 0.0957394339934444
$ws.Range("M6").Value = 0.095698985971771
$ws.Range("N6").Value = 0.0962324343265347
$ws.Range("O6").Value = 0.0970751941528084
$ws.Range("P6").Value = 0.0979732336447026
$ws.Range("Q6").Value = 0.0989242020486764
$ws.Range("R6").Value = 0.0999036965791406
$ws.Range("S6").Value = 0.100908086966793
$ws.Range("T6").Value = 0.101873951314795
$ws.Range("U6").Value = 0.102860165637175
$ws.Range("V6").Value = 0.103874281898775
$ws.Range("W6").Value = 0.104807081882312
$ws.Range("X6").Value = 0.105669997459209
$ws.Range("Y6").Value = 0.106396378514836
$ws.Range("Z6").Value = 0.10705410565214
$ws.Range("AA6").Value = 0.10762567491518
$ws.Range("AB6").Value = 0.108071655547271
$ws.Range("AC6").Value = 0.108397644061209
$ws.Range("AD6").Value = 0.584807391062041

# Remove AD4 cell (no longer present; row 4 spans 1:29 instead of 1:30)
$ws.Range("AD4").ClearContents()